{"js": "// Append \" (Changed main)\" to the end of the first paragraph\n// (\"This is a Microsoft word document.\"), split across three new runs\n// (\" (\", \"Changed main\", \")\") immediately following the existing run,\n// matching the target OOXML:\n//   <w:r><w:t>This is a Microsoft word document.</w:t></w:r>\n//   <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n//   <w:r><w:t>Changed main</w:t></w:r>\n//   <w:r><w:t>)</w:t></w:r>\n\nconst body = context.document.body;\nconst firstParagraph = body.paragraphs.getFirst();\n\n// Caret just past the existing text, before the paragraph mark.\nconst insertionPoint = firstParagraph.getRange(\"End\");\n\n// A flat-OPC WordProcessingML package that, once merged in, appends three\n// sibling runs (kept separate, rather than being coalesced into the\n// preceding run the way insertText()/InsertAfter() would).\nconst flatOpcXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n            <w:r><w:t>Changed main</w:t></w:r>\n            <w:r><w:t>)</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ninsertionPoint.insertOoxml(flatOpcXml, \"End\");\nawait context.sync();\n", "ps1": "# Append \" (Changed main)\" to the end of the first paragraph\n# (\"This is a Microsoft word document.\"), split across three new runs\n# (\" (\", \"Changed main\", \")\") immediately following the existing run,\n# matching the target OOXML:\n#   <w:r><w:t>This is a Microsoft word document.</w:t></w:r>\n#   <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n#   <w:r><w:t>Changed main</w:t></w:r>\n#   <w:r><w:t>)</w:t></w:r>\n\n$d = $word.ActiveDocument\n$firstParagraph = $d.Paragraphs(1)\n$insertionPoint = $firstParagraph.Range\n\n# A flat-OPC WordProcessingML package that, once merged in, appends three\n# sibling runs (kept separate, rather than being coalesced into the\n# preceding run the way Range.InsertAfter()/Range.Text would).\n$flatOpcXml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\"> (</w:t></w:r>\n            <w:r><w:t>Changed main</w:t></w:r>\n            <w:r><w:t>)</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$insertionPoint.InsertXML($flatOpcXml, \"End\")\n"}
